# Applies the weekly report update:
#  - Updates the "Report Generated On" timestamp in D5
#  - Zeroes out the Total Billed Amount (C8) and the per-row / total
#    pricing cells (H16, H17) to reflect the no-violation scenario

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:49 PM"

$ws.Range("C8").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
